$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.797.43'
$ws.Range("E2").Value = '  -1.85%  '

$ws.Range("D3").Value = '3.555.36'
$ws.Range("E3").Value = '  -3.55%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '190.48'
$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '567.85'
$ws.Range("E6").Value = '  -5.36%  '

$ws.Range("D7").Value = '3.548.96'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.611'
$ws.Range("E8").Value = '  -2.20%  '

$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.670'
$ws.Range("E10").Value = '  -5.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '55.38'
$ws.Range("E11").Value = '  -5.05%  '

$ws.Range("E12").Value = '  -4.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000266'
$ws.Range("E13").Value = '  -4.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.77'
$ws.Range("E14").Value = '  -4.73%  '

$ws.Range("D15").Value = '4.127.26'
$ws.Range("E15").Value = '  -3.36%  '

$ws.Range("D16").Value = '3.561.57'
$ws.Range("E16").Value = '  -3.41%  '

$ws.Range("E17").Value = '  -1.47%  '

$ws.Range("D18").Value = '66.763.34'
$ws.Range("E18").Value = '  -1.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.11'
$ws.Range("E19").Value = '  -5.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.07'
$ws.Range("E20").Value = '  -4.32%  '

$ws.Range("E21").Value = '  -6.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '397.56'
$ws.Range("E22").Value = '  -1.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.14'
$ws.Range("E23").Value = '  -7.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.49'
$ws.Range("E24").Value = '  -3.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.45'
$ws.Range("E25").Value = '  -1.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.90'
$ws.Range("E26").Value = '  -3.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.33'
$ws.Range("E27").Value = '  -2.53%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.08'
$ws.Range("E28").Value = '  +0.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.60'
$ws.Range("E29").Value = '  -3.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.68'
$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.86'
$ws.Range("E31").Value = '  -5.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '30.99'
$ws.Range("E32").Value = '  -3.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '636.26'
$ws.Range("E33").Value = '  +2.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.02'
$ws.Range("E34").Value = '  -3.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '63.68'
$ws.Range("E35").Value = '  -5.24%  '

$ws.Range("E36").Value = '  -4.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '41.85'
$ws.Range("E37").Value = '  -8.64%  '

$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("D40").Value = '0.0₃0754'
$ws.Range("E40").Value = '  -4.78%  '

$ws.Range("D41").Value = '3.186.44'
$ws.Range("E41").Value = '  +11.58%  '

$ws.Range("E42").Value = '  -2.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.65'
$ws.Range("E44").Value = '  +2.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.92'
$ws.Range("E45").Value = '  -0.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0410'
$ws.Range("E46").Value = '  -4.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.10'
$ws.Range("E47").Value = '  -4.28%  '

$ws.Range("E48").Value = '  -6.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.09'
$ws.Range("E49").Value = '  -2.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.46'
$ws.Range("E50").Value = '  -6.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.51'
$ws.Range("E51").Value = '  -5.33%  '
